$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 16
$ws.Cells.Item(16, 1).Value = ' Sharjah'
$ws.Cells.Item(16, 2).Value = ' October 26 2020'
$ws.Cells.Item(16, 3).Value = 'Kings XI won by 8 wickets (with 7 balls remaining)'
$ws.Cells.Item(16, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(16, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(16, 6).Value = 'Shubman Gill '
Set-TextCell 16 7 '57'
Set-TextCell 16 8 '45'
Set-TextCell 16 9 '3'
Set-TextCell 16 10 '4'
Set-TextCell 16 11 '126.66'

# Row 17
$ws.Cells.Item(17, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(17, 2).Value = ' September 23 2020'
$ws.Cells.Item(17, 3).Value = 'Mumbai won by 49 runs'
$ws.Cells.Item(17, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(17, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(17, 6).Value = 'Shubman Gill '
Set-TextCell 17 7 '7'
Set-TextCell 17 8 '11'
Set-TextCell 17 9 '1'
Set-TextCell 17 10 '0'
Set-TextCell 17 11 '63.63'

# Row 18
$ws.Cells.Item(18, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(18, 2).Value = ' October 18 2020'
$ws.Cells.Item(18, 3).Value = 'Match tied (KKR won the one-over eliminator)'
$ws.Cells.Item(18, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(18, 5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(18, 6).Value = 'Shubman Gill '
Set-TextCell 18 7 '36'
Set-TextCell 18 8 '37'
Set-TextCell 18 9 '5'
Set-TextCell 18 10 '0'
Set-TextCell 18 11 '97.29'

# Row 19
$ws.Cells.Item(19, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(19, 2).Value = ' October 29 2020'
$ws.Cells.Item(19, 3).Value = 'Super Kings won by 6 wickets'
$ws.Cells.Item(19, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(19, 5).Value = 'Chennai Super Kings'
$ws.Cells.Item(19, 6).Value = 'Shubman Gill '
Set-TextCell 19 7 '26'
Set-TextCell 19 8 '17'
Set-TextCell 19 9 '4'
Set-TextCell 19 10 '0'
Set-TextCell 19 11 '152.94'

# Row 20
$ws.Cells.Item(20, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(20, 2).Value = ' October 24 2020'
$ws.Cells.Item(20, 3).Value = 'KKR won by 59 runs'
$ws.Cells.Item(20, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(20, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(20, 6).Value = 'Shubman Gill '
Set-TextCell 20 7 '9'
Set-TextCell 20 8 '8'
Set-TextCell 20 9 '2'
Set-TextCell 20 10 '0'
Set-TextCell 20 11 '112.50'

# Row 21
$ws.Cells.Item(21, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(21, 2).Value = ' September 26 2020'
$ws.Cells.Item(21, 3).Value = 'KKR won by 7 wickets (with 12 balls remaining)'
$ws.Cells.Item(21, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(21, 5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(21, 6).Value = 'Shubman Gill '
Set-TextCell 21 7 '70'
Set-TextCell 21 8 '62'
Set-TextCell 21 9 '5'
Set-TextCell 21 10 '2'
Set-TextCell 21 11 '112.90'

# Row 22
$ws.Cells.Item(22, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(22, 2).Value = ' September 30 2020'
$ws.Cells.Item(22, 3).Value = 'KKR won by 37 runs'
$ws.Cells.Item(22, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(22, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(22, 6).Value = 'Shubman Gill '
Set-TextCell 22 7 '47'
Set-TextCell 22 8 '34'
Set-TextCell 22 9 '5'
Set-TextCell 22 10 '1'
Set-TextCell 22 11 '138.23'

# Row 23
$ws.Cells.Item(23, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(23, 2).Value = ' October 16 2020'
$ws.Cells.Item(23, 3).Value = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$ws.Cells.Item(23, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(23, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(23, 6).Value = 'Shubman Gill '
Set-TextCell 23 7 '21'
Set-TextCell 23 8 '23'
Set-TextCell 23 9 '2'
Set-TextCell 23 10 '0'
Set-TextCell 23 11 '91.30'

# Row 24
$ws.Cells.Item(24, 1).Value = ' Dubai (DSC)'
$ws.Cells.Item(24, 2).Value = ' November 01 2020'
$ws.Cells.Item(24, 3).Value = 'KKR won by 60 runs'
$ws.Cells.Item(24, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(24, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(24, 6).Value = 'Shubman Gill '
Set-TextCell 24 7 '36'
Set-TextCell 24 8 '24'
Set-TextCell 24 9 '6'
Set-TextCell 24 10 '0'
Set-TextCell 24 11 '150.00'

# Row 25
$ws.Cells.Item(25, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(25, 2).Value = ' October 10 2020'
$ws.Cells.Item(25, 3).Value = 'KKR won by 2 runs'
$ws.Cells.Item(25, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(25, 5).Value = 'Kings XI Punjab'
$ws.Cells.Item(25, 6).Value = 'Shubman Gill '
Set-TextCell 25 7 '57'
Set-TextCell 25 8 '47'
Set-TextCell 25 9 '5'
Set-TextCell 25 10 '0'
Set-TextCell 25 11 '121.27'

# Row 26
$ws.Cells.Item(26, 1).Value = ' Sharjah'
$ws.Cells.Item(26, 2).Value = ' October 03 2020'
$ws.Cells.Item(26, 3).Value = 'Capitals won by 18 runs'
$ws.Cells.Item(26, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(26, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(26, 6).Value = 'Shubman Gill '
Set-TextCell 26 7 '28'
Set-TextCell 26 8 '22'
Set-TextCell 26 9 '2'
Set-TextCell 26 10 '1'
Set-TextCell 26 11 '127.27'

# Row 27
$ws.Cells.Item(27, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(27, 2).Value = ' October 21 2020'
$ws.Cells.Item(27, 3).Value = 'RCB won by 8 wickets (with 39 balls remaining)'
$ws.Cells.Item(27, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(27, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(27, 6).Value = 'Shubman Gill '
Set-TextCell 27 7 '1'
Set-TextCell 27 8 '6'
Set-TextCell 27 9 '0'
Set-TextCell 27 10 '0'
Set-TextCell 27 11 '16.66'

# Row 28
$ws.Cells.Item(28, 1).Value = ' Sharjah'
$ws.Cells.Item(28, 2).Value = ' October 12 2020'
$ws.Cells.Item(28, 3).Value = 'RCB won by 82 runs'
$ws.Cells.Item(28, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(28, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(28, 6).Value = 'Shubman Gill '
Set-TextCell 28 7 '34'
Set-TextCell 28 8 '25'
Set-TextCell 28 9 '3'
Set-TextCell 28 10 '1'
Set-TextCell 28 11 '136.00'

# Row 29
$ws.Cells.Item(29, 1).Value = ' Abu Dhabi'
$ws.Cells.Item(29, 2).Value = ' October 07 2020'
$ws.Cells.Item(29, 3).Value = 'KKR won by 10 runs'
$ws.Cells.Item(29, 4).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(29, 5).Value = 'Chennai Super Kings'
$ws.Cells.Item(29, 6).Value = 'Shubman Gill '
Set-TextCell 29 7 '11'
Set-TextCell 29 8 '12'
Set-TextCell 29 9 '1'
Set-TextCell 29 10 '0'
Set-TextCell 29 11 '91.66'
